$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns per refreshed cryptos data feed.
# For Price values that parse as plain numbers, force text via a temporary
# "@" (Text) number format + restore the original cell style afterward so the
# stored type stays a string (matching the sheet's existing inlineStr cells)
# without leaving the style index changed.

# Row 2
$ws.Range("D2").Value = "27.573.30"
$ws.Range("E2").Value = "  -1.43%  "

# Row 3
$ws.Range("D3").Value = "1.629.66"
$ws.Range("E3").Value = "  -0.84%  "

# Row 4
$style = $ws.Range("D4").Style
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = $style
$ws.Range("E4").Value = "  -0.08%  "

# Row 5
$style = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.19"
$ws.Range("D5").Style = $style
$ws.Range("E5").Value = "  -0.76%  "

# Row 6
$ws.Range("E6").Value = "  -0.95%  "

# Row 7
$ws.Range("E7").Value = "  +0.03%  "

# Row 8
$style = $ws.Range("D8").Style
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "22.98"
$ws.Range("D8").Style = $style

# Row 9
$ws.Range("E9").Value = "  -0.56%  "

# Row 10
$ws.Range("E10").Value = "  -0.39%  "

# Row 11
$style = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0863"
$ws.Range("D11").Style = $style
$ws.Range("E11").Value = "  -3.14%  "

# Row 12
$ws.Range("D12").Value = "1.859.66"
$ws.Range("E12").Value = "  -0.87%  "

# Row 13
$ws.Range("D13").Value = "1.633.03"
$ws.Range("E13").Value = "  -0.43%  "

# Row 14
$ws.Range("E14").Value = "  -0.59%  "

# Row 15
$style = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.556"
$ws.Range("D15").Style = $style
$ws.Range("E15").Value = "  -0.99%  "

# Row 16
$style = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.08"
$ws.Range("D16").Style = $style
$ws.Range("E16").Value = "  +0.53%  "

# Row 17
$ws.Range("D17").Value = "27.537.33"
$ws.Range("E17").Value = "  -1.55%  "

# Row 18
$style = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "228.78"
$ws.Range("D18").Style = $style
$ws.Range("E18").Value = "  -2.07%  "

# Row 19
$ws.Range("D19").Value = "0.0₃0718"
$ws.Range("E19").Value = "  -0.74%  "

# Row 20
$ws.Range("E20").Value = "  -1.68%  "

# Row 21
$ws.Range("E21").Value = "  -0.03%  "

# Row 22
$style = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.67"
$ws.Range("D22").Style = $style
$ws.Range("E22").Value = "  +6.23%  "

# Row 23
$style = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.36"
$ws.Range("D23").Style = $style
$ws.Range("E23").Value = "  +1.06%  "

# Row 24
$ws.Range("E24").Value = "  +2.97%  "

# Row 25
$style = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "149.02"
$ws.Range("D25").Style = $style
$ws.Range("E25").Value = "  -1.03%  "

# Row 26
$ws.Range("E26").Value = "  -1.42%  "

# Row 28
$style = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.59"
$ws.Range("D28").Style = $style
$ws.Range("E28").Value = "  -0.65%  "

# Row 29
$ws.Range("E29").Value = "  -0.02%  "

# Row 30
$ws.Range("E30").Value = "  -0.81%  "

# Row 31
$ws.Range("E31").Value = "  -0.71%  "

# Row 32
$ws.Range("E32").Value = "  -1.54%  "

# Row 33
$ws.Range("D33").Value = "1.463.54"
$ws.Range("E33").Value = "  -0.54%  "

# Row 34
$ws.Range("E34").Value = "  -0.88%  "

# Row 35
$ws.Range("E35").Value = "  -0.82%  "

# Row 36
$style = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.31"
$ws.Range("D36").Style = $style
$ws.Range("E36").Value = "  -1.54%  "

# Row 37
$style = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.878"
$ws.Range("D37").Style = $style
$ws.Range("E37").Value = "  -0.53%  "

# Row 38
$ws.Range("E38").Value = "  -1.87%  "

# Row 39
$ws.Range("E39").Value = "  -0.88%  "

# Row 40
$style = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.918"
$ws.Range("D40").Style = $style
$ws.Range("E40").Value = "  -0.01%  "

# Row 41
$ws.Range("E41").Value = "  +0.71%  "

# Row 42
$ws.Range("E42").Value = "  +0.11%  "

# Row 43
$style = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "68.50"
$ws.Range("D43").Style = $style
$ws.Range("E43").Value = "  -1.82%  "

# Row 44
$ws.Range("E44").Value = "  +0.30%  "

# Row 45
$ws.Range("E45").Value = "  -1.03%  "

# Row 46
$ws.Range("E46").Value = "  -1.14%  "

# Row 47
$ws.Range("D47").Value = "1.769.15"
$ws.Range("E47").Value = "  -0.86%  "

# Row 48
$style = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.74"
$ws.Range("D48").Style = $style
$ws.Range("E48").Value = "  +1.48%  "

# Row 49
$style = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "87.30"
$ws.Range("D49").Style = $style
$ws.Range("E49").Value = "  +0.66%  "

# Row 50
$ws.Range("E50").Value = "  -0.90%  "

# Row 51
$style = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0992"
$ws.Range("D51").Style = $style
$ws.Range("E51").Value = "  -0.26%  "
